$d = $word.ActiveDocument

# --- Paragraph 1: merge "...CAT tools" + "." runs into a single run ending with a period ---
$p1 = $d.Paragraphs(1)
$p1.Range.Find.Execute("CAT tools.", $true, $false, $false, $false, $false, $true, 1, $false, "CAT tools.", 2)

# --- Paragraph 2: merge "...XTM mobile app" + "." runs into a single run ending with a period ---
$p2 = $d.Paragraphs(2)
$p2.Range.Find.Execute("XTM mobile app.", $true, $false, $false, $false, $false, $true, 1, $false, "XTM mobile app.", 2)

# --- Insert a brand-new paragraph right after paragraph 2 for the "Subcontracting..." sentence ---
$p2 = $d.Paragraphs(2)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs(3)
$p3.Range.InsertBefore("Subcontracting that includes connecting with translation vendors, rate cards, cost estimates and purchase orders.")

# --- Remove the hidden _GoBack bookmark that used to sit at the end of the old "Subcontracting" paragraph ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- The old third paragraph (now the fourth) becomes the new closing line ---
$p4 = $d.Paragraphs(4)
$p4.Range.Find.Execute("Subcontracting that includes connecting with translation vendors, rate cards, cost estimates and purchase orders.", $true, $false, $false, $false, $false, $true, 1, $false, "Pretty cool, right?", 2)
